$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2687
$ws.Range("F9").Value = 261
$ws.Range("F10").Value = 5981
$ws.Range("F13").Value = 4909
$ws.Range("F15").Value = 92
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 2520
$ws.Range("F18").Value = 1316
$ws.Range("F19").Value = 490
$ws.Range("F20").Value = 1198
$ws.Range("F22").Value = 271
$ws.Range("F24").Value = 123
$ws.Range("F26").Value = 215
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 1336
$ws.Range("F31").Value = 5
$ws.Range("F32").Value = 2067
$ws.Range("F33").Value = 280
$ws.Range("F34").Value = 560
$ws.Range("F35").Value = 52
$ws.Range("F37").Value = 1454
$ws.Range("F38").Value = 605
$ws.Range("F40").Value = 544
$ws.Range("F41").Value = 229
$ws.Range("F42").Value = 1710
$ws.Range("F43").Value = 2496
$ws.Range("F44").Value = 52
$ws.Range("F45").Value = 102
$ws.Range("F48").Value = 63

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 5
$ws.Range("F11").Value = 196
$ws.Range("F15").Value = 146
$ws.Range("F23").Value = 324
$ws.Range("F31").Value = 3
$ws.Range("F35").Value = 11
$ws.Range("F38").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1687
$ws.Range("F8").Value = 1426
$ws.Range("F9").Value = 1795
$ws.Range("F10").Value = 2381
$ws.Range("F11").Value = 785
$ws.Range("F12").Value = 668

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1687
$ws.Range("F7").Value = 2687
$ws.Range("F9").Value = 1426
$ws.Range("F10").Value = 2381
$ws.Range("F11").Value = 5981
$ws.Range("F12").Value = 785
$ws.Range("F16").Value = 4909
$ws.Range("F17").Value = 92
$ws.Range("F18").Value = 2520
$ws.Range("F19").Value = 1316
$ws.Range("F20").Value = 490
$ws.Range("F21").Value = 1198
$ws.Range("F22").Value = 271
$ws.Range("F24").Value = 123
$ws.Range("F26").Value = 215
$ws.Range("F28").Value = 1336
$ws.Range("F29").Value = 2067
$ws.Range("F30").Value = 280
$ws.Range("F31").Value = 560
$ws.Range("F34").Value = 1454
$ws.Range("F35").Value = 605
$ws.Range("F37").Value = 544
$ws.Range("F40").Value = 229
$ws.Range("F42").Value = 1710
$ws.Range("F43").Value = 2496
$ws.Range("F44").Value = 102
$ws.Range("F47").Value = 63
$ws.Range("F49").Value = 11

